$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3915976666666667
$ws.Range("H2").Value = 1.174793
$ws.Range("I2").Value = 0.02606065131430495
$ws.Range("J2").Value = 0.02606065131430495
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 0.5650903137113333
$ws.Range("R2").Value = 5.085812823402001
$ws.Range("S2").Value = 0.0007548907580061876
$ws.Range("T2").Value = 0.0007548907580061875

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3915976666666667
$ws.Range("H3").Value = 1.174793
$ws.Range("I3").Value = 0.02606065131430495
$ws.Range("J3").Value = 0.02606065131430495
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 11.43602807139989
$ws.Range("R3").Value = 102.924252642599
$ws.Range("S3").Value = 0.0152771188780438
$ws.Range("T3").Value = 0.0152771188780438

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3915976666666667
$ws.Range("H4").Value = 1.174793
$ws.Range("I4").Value = 0.02606065131430495
$ws.Range("J4").Value = 0.02606065131430495
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 7.507163403393001
$ws.Range("R4").Value = 67.56447063053702
$ws.Range("S4").Value = 0.01002864167825497
$ws.Range("T4").Value = 0.01002864167825497

# Row 5
$ws.Range("I5").Value = 0.4187506438669658
$ws.Range("J5").Value = 0.4187506438669658
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 9.080046766893998
$ws.Range("R5").Value = 81.720420902046
$ws.Range("S5").Value = 0.01212981928777799
$ws.Range("T5").Value = 0.01212981928777798

# Row 6
$ws.Range("I6").Value = 0.4187506438669658
$ws.Range("J6").Value = 0.4187506438669658
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("S6").Value = 0.2454774936151145
$ws.Range("T6").Value = 0.2454774936151145

# Row 7
$ws.Range("I7").Value = 0.4187506438669658
$ws.Range("J7").Value = 0.4187506438669658
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 120.627434474739
$ws.Range("R7").Value = 1085.646910272651
$ws.Range("S7").Value = 0.1611433309640734
$ws.Range("T7").Value = 0.1611433309640734

# Row 8
$ws.Range("G8").Value = 8.342485333333334
$ws.Range("H8").Value = 25.027456
$ws.Range("I8").Value = 0.5551887048187292
$ws.Range("J8").Value = 0.5551887048187292
$ws.Range("M8").Value = 1.443038
$ws.Range("N8").Value = 4.329114
$ws.Range("O8").Value = 0.0289666880885598
$ws.Range("P8").Value = 0.0289666880885598
$ws.Range("Q8").Value = 12.03852335044267
$ws.Range("R8").Value = 108.346710153984
$ws.Range("S8").Value = 0.01608197804277562
$ws.Range("T8").Value = 0.01608197804277562

# Row 9
$ws.Range("G9").Value = 8.342485333333334
$ws.Range("H9").Value = 25.027456
$ws.Range("I9").Value = 0.5551887048187292
$ws.Range("J9").Value = 0.5551887048187292
$ws.Range("N9").Value = 87.61054300000001
$ws.Range("O9").Value = 0.5862140087672342
$ws.Range("P9").Value = 0.5862140087672342
$ws.Range("Q9").Value = 243.6298900076231
$ws.Range("R9").Value = 2192.669010068608
$ws.Range("S9").Value = 0.3254593962740759
$ws.Range("T9").Value = 0.3254593962740759

# Row 10
$ws.Range("G10").Value = 8.342485333333334
$ws.Range("H10").Value = 25.027456
$ws.Range("I10").Value = 0.5551887048187292
$ws.Range("J10").Value = 0.5551887048187292
$ws.Range("M10").Value = 19.170603
$ws.Range("N10").Value = 57.511809
$ws.Range("O10").Value = 0.384819303144206
$ws.Range("P10").Value = 0.384819303144206
$ws.Range("Q10").Value = 159.930474358656
$ws.Range("R10").Value = 1439.374269227904
$ws.Range("S10").Value = 0.2136473305018777
$ws.Range("T10").Value = 0.2136473305018776
